$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F holds "想去人数" (want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5536
$ws1.Range("F3").Value = 610
$ws1.Range("F4").Value = 12370
$ws1.Range("F5").Value = 301
$ws1.Range("F7").Value = 187
$ws1.Range("F8").Value = 355
$ws1.Range("F9").Value = 1140
$ws1.Range("F10").Value = 108

# Sheet "全部类型" (All types) - same events, same column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5536
$ws4.Range("F4").Value = 610
$ws4.Range("F6").Value = 12370
$ws4.Range("F7").Value = 301
$ws4.Range("F9").Value = 187
$ws4.Range("F12").Value = 355
$ws4.Range("F13").Value = 1140
$ws4.Range("F15").Value = 108
